$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C11").Value = "22:20-22:25"
$ws.Range("C12").Value = "22:25-22:30"
